# Applies the diff to the Slovenia 2-SNL 2023-2024 workbook:
#  1) Rows 5 and 6 have their match data (columns F:V) swapped.
#  2) Rows 26 and 27 have their match data (columns F:V) swapped.
#  3) Four new match rows (62-65) are appended after the existing last row (61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($rowA, $rowB, $firstCol, $lastCol)

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $valA = $ws.Cells.Item($rowA, $col).Value2
        $valB = $ws.Cells.Item($rowB, $col).Value2
        $ws.Cells.Item($rowA, $col).Value2 = $valB
        $ws.Cells.Item($rowB, $col).Value2 = $valA
    }
}

# --- 1) swap rows 5 & 6, columns F (6) through V (22) ---
Swap-RowData 5 6 6 22

# --- 2) swap rows 26 & 27, columns F (6) through V (22) ---
Swap-RowData 26 27 6 22

# --- 3) append 4 new rows (62-65) after row 61, copying row 61's formatting ---
$ws.Range("A61:V61").Copy()
$ws.Range("A62:V65").PasteSpecial(-4122)

$newRows = @(
    @{ Row=62; A=61; E=45190.66666666666; F="Bilje";          G=1; H="Triglav";  I=1;
       J=2.38; K="20/09/2023 03:12"; L=2.69; M="21/09/2023 15:59";
       N=3.24; O="20/09/2023 03:12"; P=3.06; Q="21/09/2023 15:59";
       R=2.51; S="20/09/2023 03:12"; T=2.6;  U="21/09/2023 15:59";
       V="https://www.betexplorer.com/football/slovenia/2-snl/bilje-triglav/OKlbEQ09/" },

    @{ Row=63; A=62; E=45190.66666666666; F="Jadran Dekani";  G=0; H="Beltinci"; I=3;
       J=2.74; K="20/09/2023 03:12"; L=3.03; M="21/09/2023 15:59";
       N=3.07; O="20/09/2023 03:12"; P=3.18; Q="21/09/2023 15:59";
       R=2.3;  S="20/09/2023 03:12"; T=2.27; U="21/09/2023 14:40";
       V="https://www.betexplorer.com/football/slovenia/2-snl/jadran-dekani-beltinci/lSoABOpS/" },

    @{ Row=64; A=63; E=45190.66666666666; F="Nafta";          G=3; H="Tolmin";   I=0;
       J=1.32; K="20/09/2023 03:12"; L=1.47; M="21/09/2023 15:58";
       N=4.8;  O="20/09/2023 03:12"; P=4.17; Q="21/09/2023 15:59";
       R=6.17; S="20/09/2023 03:12"; T=6.11; U="21/09/2023 15:59";
       V="https://www.betexplorer.com/football/slovenia/2-snl/nafta-tolmin/Aym6C4VL/" },

    @{ Row=65; A=64; E=45190.85416666666; F="Rudar";          G=2; H="Grosuplje"; I=0;
       J=3.12; K="20/09/2023 07:42"; L=4.23; M="21/09/2023 20:26";
       N=3.33; O="20/09/2023 07:42"; P=3.51; Q="21/09/2023 20:30";
       R=2.01; S="20/09/2023 07:42"; T=1.78; U="21/09/2023 20:30";
       V="https://www.betexplorer.com/football/slovenia/2-snl/rudar-grosuplje/2Bm2DpGF/" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $r.A
    $ws.Cells.Item($row, 2).Value2  = "slovenia"
    $ws.Cells.Item($row, 3).Value2  = "2-snl"
    $ws.Cells.Item($row, 4).Value2  = "2023-2024"
    $ws.Cells.Item($row, 5).Value2  = $r.E
    $ws.Cells.Item($row, 6).Value2  = $r.F
    $ws.Cells.Item($row, 7).Value2  = $r.G
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $r.N
    $ws.Cells.Item($row, 15).Value2 = $r.O
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $r.Q
    $ws.Cells.Item($row, 18).Value2 = $r.R
    $ws.Cells.Item($row, 19).Value2 = $r.S
    $ws.Cells.Item($row, 20).Value2 = $r.T
    $ws.Cells.Item($row, 21).Value2 = $r.U
    $ws.Cells.Item($row, 22).Value2 = $r.V
}
